$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the bold/bordered header style (row 1) ---
# Revert header row formatting back to the default/Normal style so the
# extra font/border/cellXf entries used only by the header go unused.
$ws.Range("A1:R1").Style = "Normal"

# --- Clear the "Unnamed: 0" header label in A1 ---
$ws.Range("A1").ClearContents()

# --- Corrected data-cleaning values (pre/post/total fixation metrics) ---
# Row 3 (Revisit count)
$ws.Range("C3").Value = 39
$ws.Range("D3").Value = 46
$ws.Range("K3").Value = 12
$ws.Range("N3").Value = 54

# Row 4 (Fixation count)
$ws.Range("C4").Value = 241
$ws.Range("D4").Value = 145
$ws.Range("K4").Value = 15
$ws.Range("N4").Value = 351

# Row 6 (Dwell time (ms))
$ws.Range("C6").Value = 68842.99000000001
$ws.Range("D6").Value = 39078.91
$ws.Range("K6").Value = 4220.13
$ws.Range("N6").Value = 126217.96

# Row 7 (Dwell time (%))
$ws.Range("B7").Value = 1.6
$ws.Range("C7").Value = 21.73
$ws.Range("D7").Value = 12.34
$ws.Range("E7").Value = 7.85
$ws.Range("G7").Value = 2.91
$ws.Range("H7").Value = 1.04
$ws.Range("I7").Value = 3.65
$ws.Range("J7").Value = 2.65
$ws.Range("K7").Value = 1.33
$ws.Range("L7").Value = 4.16
$ws.Range("M7").Value = 1.08
$ws.Range("N7").Value = 39.84
$ws.Range("O7").Value = 0.4
$ws.Range("P7").Value = 0.13
$ws.Range("Q7").Value = 0.21
$ws.Range("R7").Value = 0.52

# Row 8 (Fixation duration (ms))
$ws.Range("C8").Value = 285.66
$ws.Range("D8").Value = 269.51
$ws.Range("K8").Value = 281.34
$ws.Range("N8").Value = 359.6

# --- Drop the trailing blank rows 11-13 ---
$ws.Rows("11:13").Delete()
